# Append a new portfolio-update row (2025-09-19) to the bottom of the
# sheet's data table, mirroring the existing rows: column A holds the
# date as literal text (not an auto-converted date serial), columns
# B-D hold plain numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.UsedRange.Row + $ws.UsedRange.Rows.Count, 1).Row

# Force column A to be treated as text so "2025-09-19" isn't coerced into
# a date serial number, then restore the default "Normal" style so the
# cell doesn't end up with an explicit/custom number-format style (matches
# the unstyled look of every other data row).
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025-09-19"
$ws.Range("A" + $newRow).Style = "Normal"

$ws.Range("B" + $newRow).Value = 60.40000152587891
$ws.Range("C" + $newRow).Value = 707.4500122070312
$ws.Range("D" + $newRow).Value = 336.5499877929688
